# Adjust input tab naming to match strategic model
#
# - Rename "InitialDisposalCapacity" -> "DisposalCapacity"
# - Remove the (effectively empty) "DriveDistances" placeholder tab
# - Rename "DriveTimes" -> "TruckingTime"
# - Rename "PipingOperationalCost" -> "PipelineOperationalCost"
# - Reorder the operational-cost / trucking tabs so that costs are grouped
#   together, followed by the trucking-hours / trucking-time tabs
# - Give "CompletionsDemand" the same tab color as its neighboring input tabs

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Renames -------------------------------------------------------------
$wb.Worksheets.Item("InitialDisposalCapacity").Name = "DisposalCapacity"
$wb.Worksheets.Item("PipingOperationalCost").Name = "PipelineOperationalCost"
$wb.Worksheets.Item("DriveTimes").Name = "TruckingTime"

# --- Delete the obsolete DriveDistances tab ------------------------------
[void]$wb.Worksheets.Item("DriveDistances").Delete()

# --- Reorder tabs ---------------------------------------------------------
# Walk the target order, moving each sheet (in turn) to right after the
# previously-placed sheet, anchored on the sheet that doesn't move. Re-fetch
# the anchor worksheet by name each time (rather than keeping the old object
# reference around) so each Move() sees the worksheet's current position.
$order = @(
    "TreatmentOperationalCost",
    "ReuseOperationalCost",
    "PipelineOperationalCost",
    "FreshSourcingCost",
    "PadStorageCost",
    "TruckingHourlyCost",
    "TruckingTime"
)

$anchorName = "DisposalOperationalCost"
foreach ($name in $order) {
    $ws = $wb.Worksheets.Item($name)
    $anchor = $wb.Worksheets.Item($anchorName)
    $ws.Move($null, $anchor)
    $anchorName = $name
}

# --- Match CompletionsDemand tab color to its sibling input tabs ---------
# (ProductionRates/PadRates/FlowbackRates use theme 7, tint 0.8 ~ RGB FFF2CC)
$wb.Worksheets.Item("CompletionsDemand").Tab.Color = 13431551
